$wb = $excel.ActiveWorkbook

# Worksheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3025
$ws.Range("F4").Value = 2348
$ws.Range("F7").Value = 874
$ws.Range("F8").Value = 146
$ws.Range("F11").Value = 2697
$ws.Range("F13").Value = 1575
$ws.Range("F14").Value = 7258
$ws.Range("F16").Value = 7410
$ws.Range("F17").Value = 10
$ws.Range("F18").Value = 26
$ws.Range("F19").Value = 5787
$ws.Range("F20").Value = 3167
$ws.Range("F21").Value = 3539
$ws.Range("F22").Value = 10
$ws.Range("F24").Value = 254
$ws.Range("F26").Value = 1992
$ws.Range("F27").Value = 89
$ws.Range("F30").Value = 242
$ws.Range("F31").Value = 715
$ws.Range("F32").Value = 50
$ws.Range("F33").Value = 2512
$ws.Range("F34").Value = 1328
$ws.Range("F35").Value = 2990
$ws.Range("F36").Value = 98
$ws.Range("F38").Value = 186
$ws.Range("F39").Value = 436
$ws.Range("F40").Value = 1160
$ws.Range("F43").Value = 560

# Worksheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 98
$ws.Range("F11").Value = 383

# Worksheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 3025
$ws.Range("F5").Value = 2348
$ws.Range("F8").Value = 874
$ws.Range("F9").Value = 146
$ws.Range("F14").Value = 2697
$ws.Range("F15").Value = 1575
$ws.Range("F17").Value = 7258
$ws.Range("F19").Value = 7410
$ws.Range("F20").Value = 10
$ws.Range("F21").Value = 26
$ws.Range("F22").Value = 5787
$ws.Range("F23").Value = 3167
$ws.Range("F24").Value = 3539
$ws.Range("F25").Value = 10
$ws.Range("F28").Value = 254
$ws.Range("F30").Value = 1992
$ws.Range("F35").Value = 242
$ws.Range("F36").Value = 715
$ws.Range("F37").Value = 50
$ws.Range("F38").Value = 2512
$ws.Range("F39").Value = 1328
$ws.Range("F41").Value = 2991
$ws.Range("F42").Value = 98
$ws.Range("F44").Value = 186
$ws.Range("F46").Value = 436
$ws.Range("F47").Value = 1160
$ws.Range("F49").Value = 560
